$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the ipaddress value (B2) with the new world/host IP
$ws.Range("B2").Value = "192.168.1.131"

# Update the selected cell to match the new view state
$ws.Range("B5").Select()
